$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-23: update the changed-date value
$ws.Range("C2:C23").Value = 45174
